$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5779.2
$ws.Range("J112").Value = 6376.8887
$ws.Range("L112").Value = 19130.6661
$ws.Range("N112").Value = -21346.6661
$ws.Range("H125").Value = 739.75
$ws.Range("I125").Value = 383.6
$ws.Range("J125").Value = 1333.3334
$ws.Range("K125").Value = 3452.4
$ws.Range("L125").Value = 12000.0006
$ws.Range("M125").Value = -992.4000000000001
$ws.Range("N125").Value = -16920.0006
$ws.Range("H129").Value = 913.04
$ws.Range("J129").Value = 889.19147
$ws.Range("L129").Value = 2667.57441
$ws.Range("N129").Value = -12667.57441
$ws.Range("H132").Value = 1094.3429
$ws.Range("I132").Value = 1040.7188
$ws.Range("J132").Value = 1666.3334
$ws.Range("K132").Value = 3122.1564
$ws.Range("L132").Value = 4999.0002
$ws.Range("M132").Value = -592.1564000000003
$ws.Range("N132").Value = -10059.0002
$ws.Range("H135").Value = 527.125
$ws.Range("I135").Value = 514.1429000000001
$ws.Range("K135").Value = 4627.2861
$ws.Range("M135").Value = -2092.2861
$ws.Range("H138").Value = 2867.8113
$ws.Range("J138").Value = 2805.3572
$ws.Range("L138").Value = 8416.071599999999
$ws.Range("N138").Value = -18696.0716

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4287.5
$ws.Range("I61").Value = 3099.0667
$ws.Range("K61").Value = 3099.0667
$ws.Range("M61").Value = -2887.0667
$ws.Range("H74").Value = 1120.0416
$ws.Range("I74").Value = 754.825
$ws.Range("J74").Value = 2946.125
$ws.Range("K74").Value = 754.825
$ws.Range("L74").Value = 2946.125
$ws.Range("M74").Value = 119.175
$ws.Range("N74").Value = -4694.125
$ws.Range("H77").Value = 1120.0416
$ws.Range("I77").Value = 754.825
$ws.Range("J77").Value = 2946.125
$ws.Range("K77").Value = 3774.125
$ws.Range("L77").Value = 14730.625
$ws.Range("M77").Value = 593.875
$ws.Range("N77").Value = -23466.625
$ws.Range("H88").Value = 4744.222
$ws.Range("I88").Value = 2249.5
$ws.Range("K88").Value = 2249.5
$ws.Range("M88").Value = -1843.5
$ws.Range("H91").Value = 4744.222
$ws.Range("I91").Value = 2249.5
$ws.Range("K91").Value = 2249.5
$ws.Range("M91").Value = -845.5
$ws.Range("H122").Value = 1228
$ws.Range("I122").Value = 1228
$ws.Range("K122").Value = 3684
$ws.Range("M122").Value = -1234
$ws.Range("H132").Value = 1547.1666
$ws.Range("I132").Value = 955.3
$ws.Range("K132").Value = 2865.9
$ws.Range("M132").Value = -335.8999999999996
$ws.Range("H136").Value = 4287.5
$ws.Range("I136").Value = 3099.0667
$ws.Range("K136").Value = 9297.2001
$ws.Range("M136").Value = -6747.2001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 8017.846
$ws.Range("I80").Value = 48
$ws.Range("J80").Value = 9466.909
$ws.Range("K80").Value = 48
$ws.Range("L80").Value = 9466.909
$ws.Range("M80").Value = 950
$ws.Range("N80").Value = -11462.909
$ws.Range("H83").Value = 8017.846
$ws.Range("I83").Value = 48
$ws.Range("J83").Value = 9466.909
$ws.Range("K83").Value = 240
$ws.Range("L83").Value = 47334.545
$ws.Range("M83").Value = 4752
$ws.Range("N83").Value = -57318.545
$ws.Range("H86").Value = 113457.664
$ws.Range("I86").Value = 2228.25
$ws.Range("K86").Value = 2228.25
$ws.Range("M86").Value = -1105.25
$ws.Range("H89").Value = 113457.664
$ws.Range("I89").Value = 2228.25
$ws.Range("K89").Value = 11141.25
$ws.Range("M89").Value = -5525.25

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1965
$ws.Range("I132").Value = 1103.5238
$ws.Range("K132").Value = 3310.5714
$ws.Range("M132").Value = -780.5713999999998

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 566
$ws.Range("I68").Value = 566
$ws.Range("K68").Value = 1698
$ws.Range("M68").Value = -887
$ws.Range("H71").Value = 566
$ws.Range("I71").Value = 566
$ws.Range("K71").Value = 5094
$ws.Range("M71").Value = -1038
$ws.Range("H113").Value = 10909.1
$ws.Range("J113").Value = 1070
$ws.Range("L113").Value = 3210
$ws.Range("N113").Value = -7550
$ws.Range("H131").Value = 9007.959000000001
$ws.Range("I131").Value = 532.6
$ws.Range("J131").Value = 9982.138000000001
$ws.Range("K131").Value = 1597.8
$ws.Range("L131").Value = 29946.414
$ws.Range("M131").Value = 3442.2
$ws.Range("N131").Value = -40026.414
$ws.Range("H139").Value = 15857.143
$ws.Range("I139").Value = 34000
$ws.Range("K139").Value = 102000
$ws.Range("M139").Value = -96860
$ws.Range("H140").Value = 1680.069
$ws.Range("I140").Value = 975.53845
$ws.Range("J140").Value = 2252.5
$ws.Range("K140").Value = 2926.61535
$ws.Range("L140").Value = 6757.5
$ws.Range("M140").Value = 2253.38465
$ws.Range("N140").Value = -17117.5

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17283.334
$ws.Range("I70").Value = 23825
$ws.Range("J70").Value = 4200
$ws.Range("K70").Value = 23825
$ws.Range("L70").Value = 4200
$ws.Range("M70").Value = -23555
$ws.Range("N70").Value = -4740
$ws.Range("H73").Value = 17283.334
$ws.Range("I73").Value = 23825
$ws.Range("J73").Value = 4200
$ws.Range("K73").Value = 23825
$ws.Range("L73").Value = 4200
$ws.Range("M73").Value = -22889
$ws.Range("N73").Value = -6072

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5412.625
$ws.Range("I7").Value = 3826
$ws.Range("J7").Value = 6999.25
$ws.Range("K7").Value = 3826
$ws.Range("L7").Value = 6999.25
$ws.Range("M7").Value = -3714
$ws.Range("N7").Value = -7223.25
$ws.Range("H22").Value = 2718.25
$ws.Range("J22").Value = 1907.6666
$ws.Range("L22").Value = 1907.6666
$ws.Range("N22").Value = -2497.6666
$ws.Range("H27").Value = 2718.25
$ws.Range("J27").Value = 1907.6666
$ws.Range("L27").Value = 1907.6666
$ws.Range("N27").Value = -2121.6666
$ws.Range("H46").Value = 1582.8462
$ws.Range("I46").Value = 1090
$ws.Range("J46").Value = 1672.4546
$ws.Range("K46").Value = 1090
$ws.Range("L46").Value = 1672.4546
$ws.Range("M46").Value = -902
$ws.Range("N46").Value = -2048.4546
$ws.Range("H55").Value = 309.5484
$ws.Range("I55").Value = 276.38095
$ws.Range("J55").Value = 379.2
$ws.Range("K55").Value = 276.38095
$ws.Range("L55").Value = 379.2
$ws.Range("M55").Value = -103.38095
$ws.Range("N55").Value = -725.2
$ws.Range("H68").Value = 2248.5625
$ws.Range("I68").Value = 1998.3572
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 1998.3572
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1249.3572
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 2248.5625
$ws.Range("I71").Value = 1998.3572
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 9991.786
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -6247.786
$ws.Range("N71").Value = -27488
$ws.Range("H126").Value = 5412.625
$ws.Range("I126").Value = 3826
$ws.Range("J126").Value = 6999.25
$ws.Range("K126").Value = 11478
$ws.Range("L126").Value = 20997.75
$ws.Range("M126").Value = -9008
$ws.Range("N126").Value = -25937.75

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3241.5625
$ws.Range("J132").Value = 4398.2
$ws.Range("L132").Value = 13194.6
$ws.Range("N132").Value = -18254.6

Write-Output "done"